$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45203 -> 45204, i.e. 2023-10-04 -> 2023-10-05) for every data row.
$ws.Range("C2:C45").Value = 45204
